$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume cells stay text (many values look numeric, e.g. "0.9970")
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.980.28"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "1.728.08"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").Value = "0.9970"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "239.96"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "0.9971"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4780"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("D8").Value = "0.2571"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "0.06135"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "1.719.57"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "15.92"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").Value = "0.06880"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "0.5987"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "4.414"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "76.53"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "26.936.09"
$ws.Range("E17").Value = "  +2.62%  "
$ws.Range("D18").Value = "0.9968"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "0.000007014"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").Value = "11.31"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("D21").Value = "1.938.28"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "4.373"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "8.362"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "5.073"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "140.71"
$ws.Range("D26").Value = "15.19"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "1.802"
$ws.Range("E27").Value = "  +4.07%  "
$ws.Range("D28").Value = "1.401"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "106.27"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "3.951"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").Value = "0.07900"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "3.649"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").Value = "0.04658"
$ws.Range("E33").Value = "  +4.97%  "
$ws.Range("D34").Value = "2.594"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "1.002"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").Value = "0.6123"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "0.9171"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("D38").Value = "2.519"
$ws.Range("E38").Value = "  +6.21%  "
$ws.Range("D39").Value = "1.991"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "0.9976"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "5.662"
$ws.Range("E41").Value = "  +5.67%  "
$ws.Range("D42").Value = "0.01482"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").Value = "99.27"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "0.3796"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "6.748"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "0.1145"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").Value = "0.05340"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "7.779"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").Value = "29.77"
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("D50").Value = "1.234"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").Value = "50.72"
$ws.Range("E51").Value = "  -0.82%  "
